$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 989.5625
$ws.Range("J40").Value = 799.625
$ws.Range("L40").Value = 799.625
$ws.Range("N40").Value = -1149.625

$ws.Range("H70").Value = 6594.4
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 6594.4
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 19783.2
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -20323.2

$ws.Range("H73").Value = 6594.4
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 6594.4
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 19783.2
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -21655.2

$ws.Range("H86").Value = 411227.34
$ws.Range("J86").Value = 774.5
$ws.Range("L86").Value = 774.5
$ws.Range("N86").Value = -3020.5

$ws.Range("H89").Value = 411227.34
$ws.Range("J89").Value = 774.5
$ws.Range("L89").Value = 3872.5
$ws.Range("N89").Value = -15104.5

$ws.Range("H131").Value = 1997.4762
$ws.Range("J131").Value = 3070.2727
$ws.Range("L131").Value = 9210.8181
$ws.Range("N131").Value = -19290.8181

$ws.Range("H132").Value = 1024.4865
$ws.Range("I132").Value = 963.41174
$ws.Range("J132").Value = 1716.6666
$ws.Range("K132").Value = 2890.23522
$ws.Range("L132").Value = 5149.9998
$ws.Range("M132").Value = -360.23522
$ws.Range("N132").Value = -10209.9998

$ws.Range("H137").Value = 2464.5881
$ws.Range("I137").Value = 1390.4286
$ws.Range("J137").Value = 3216.5
$ws.Range("K137").Value = 4171.2858
$ws.Range("L137").Value = 9649.5
$ws.Range("M137").Value = -1621.2858
$ws.Range("N137").Value = -14749.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3087.7817
$ws.Range("I32").Value = 2840.3044
$ws.Range("K32").Value = 2840.3044
$ws.Range("M32").Value = -2553.3044

$ws.Range("H45").Value = 1853.8182
$ws.Range("I45").Value = 1733.3334
$ws.Range("J45").Value = 1899
$ws.Range("K45").Value = 1733.3334
$ws.Range("L45").Value = 1899
$ws.Range("M45").Value = -1356.3334
$ws.Range("N45").Value = -2653

$ws.Range("H61").Value = 10861
$ws.Range("I61").Value = 8000
$ws.Range("J61").Value = 11814.667
$ws.Range("K61").Value = 8000
$ws.Range("L61").Value = 11814.667
$ws.Range("M61").Value = -7788
$ws.Range("N61").Value = -12238.667

$ws.Range("H101").Value = 45201.75
$ws.Range("J101").Value = 45201.75
$ws.Range("L101").Value = 45201.75
$ws.Range("N101").Value = -51691.75

$ws.Range("H132").Value = 2418.3462
$ws.Range("I132").Value = 1445.0588
$ws.Range("K132").Value = 4335.1764
$ws.Range("M132").Value = -1805.1764

$ws.Range("H136").Value = 10861
$ws.Range("I136").Value = 8000
$ws.Range("J136").Value = 11814.667
$ws.Range("K136").Value = 24000
$ws.Range("L136").Value = 35444.001
$ws.Range("M136").Value = -21450
$ws.Range("N136").Value = -40544.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2192.2068
$ws.Range("I20").Value = 2135.9583
$ws.Range("K20").Value = 2135.9583
$ws.Range("M20").Value = -1888.9583

$ws.Range("H94").Value = 928.8889
$ws.Range("I94").Value = 708.5714
$ws.Range("J94").Value = 1700
$ws.Range("K94").Value = 708.5714
$ws.Range("L94").Value = 1700
$ws.Range("M94").Value = -257.5714
$ws.Range("N94").Value = -2602

$ws.Range("H134").Value = 8969.096
$ws.Range("I134").Value = 9938.412
$ws.Range("K134").Value = 29815.236
$ws.Range("M134").Value = -27280.236

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1791.3871
$ws.Range("I31").Value = 1953.45
$ws.Range("K31").Value = 1953.45
$ws.Range("M31").Value = -1658.45

$ws.Range("H34").Value = 1791.3871
$ws.Range("I34").Value = 1953.45
$ws.Range("K34").Value = 1953.45
$ws.Range("M34").Value = -1751.45

$ws.Range("H58").Value = 3954305.2
$ws.Range("I58").Value = 3954305.2
$ws.Range("K58").Value = 3954305.2
$ws.Range("M58").Value = -3954102.2

$ws.Range("H62").Value = 2039.6
$ws.Range("I62").Value = 2399.75
$ws.Range("J62").Value = 599
$ws.Range("K62").Value = 2399.75
$ws.Range("L62").Value = 599
$ws.Range("M62").Value = -1775.75
$ws.Range("N62").Value = -1847

$ws.Range("H65").Value = 2039.6
$ws.Range("I65").Value = 2399.75
$ws.Range("J65").Value = 599
$ws.Range("K65").Value = 11998.75
$ws.Range("L65").Value = 2995
$ws.Range("M65").Value = -8878.75
$ws.Range("N65").Value = -9235

$ws.Range("H105").Value = 1252.1305
$ws.Range("I105").Value = 1046.625
$ws.Range("J105").Value = 1721.8572
$ws.Range("K105").Value = 1046.625
$ws.Range("L105").Value = 1721.8572
$ws.Range("M105").Value = 700.375
$ws.Range("N105").Value = -5215.8572

$ws.Range("H132").Value = 2107.5557
$ws.Range("I132").Value = 1264.0834
$ws.Range("K132").Value = 3792.2502
$ws.Range("M132").Value = -1262.2502

$ws.Range("H134").Value = 3193.5
$ws.Range("I134").Value = 3089.889
$ws.Range("J134").Value = 3380
$ws.Range("K134").Value = 9269.667000000001
$ws.Range("L134").Value = 10140
$ws.Range("M134").Value = -6734.667000000001
$ws.Range("N134").Value = -15210

$ws.Range("H136").Value = 3954305.2
$ws.Range("I136").Value = 3954305.2
$ws.Range("K136").Value = 11862915.6
$ws.Range("M136").Value = -11860365.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 125088650
$ws.Range("I50").Value = 696969
$ws.Range("J50").Value = 142858880
$ws.Range("K50").Value = 2090907
$ws.Range("L50").Value = 428576640
$ws.Range("M50").Value = -2090426
$ws.Range("N50").Value = -428577602

$ws.Range("H52").Value = 666.5
$ws.Range("J52").Value = 666.5
$ws.Range("L52").Value = 1999.5
$ws.Range("N52").Value = -2531.5

$ws.Range("H53").Value = 125088650
$ws.Range("I53").Value = 696969
$ws.Range("J53").Value = 142858880
$ws.Range("K53").Value = 2090907
$ws.Range("L53").Value = 428576640
$ws.Range("M53").Value = -2090426
$ws.Range("N53").Value = -428577602

$ws.Range("H55").Value = 13855.889
$ws.Range("I55").Value = 100004
$ws.Range("J55").Value = 3087.375
$ws.Range("K55").Value = 300012
$ws.Range("L55").Value = 9262.125
$ws.Range("M55").Value = -299835
$ws.Range("N55").Value = -9616.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4470.7144
$ws.Range("I70").Value = 4579
$ws.Range("K70").Value = 4579
$ws.Range("M70").Value = -4309

$ws.Range("H73").Value = 4470.7144
$ws.Range("I73").Value = 4579
$ws.Range("K73").Value = 4579
$ws.Range("M73").Value = -3643

$ws.Range("H97").Value = 1775.4286
$ws.Range("I97").Value = 1571.5834
$ws.Range("K97").Value = 1571.5834
$ws.Range("M97").Value = -1075.5834

$ws.Range("H122").Value = 1397.2
$ws.Range("I122").Value = 997.8570999999999
$ws.Range("J122").Value = 2329
$ws.Range("K122").Value = 2993.5713
$ws.Range("L122").Value = 6987
$ws.Range("M122").Value = -543.5712999999996
$ws.Range("N122").Value = -11887

$ws.Range("H132").Value = 1481975.8
$ws.Range("I132").Value = 2138651.2
$ws.Range("K132").Value = 6415953.600000001
$ws.Range("M132").Value = -6413423.600000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1099.6154
$ws.Range("I82").Value = 854.1818
$ws.Range("K82").Value = 854.1818
$ws.Range("M82").Value = -493.1818

$ws.Range("H85").Value = 1099.6154
$ws.Range("I85").Value = 854.1818
$ws.Range("K85").Value = 854.1818
$ws.Range("M85").Value = 393.8182

$ws.Range("H122").Value = 2958.0625
$ws.Range("I122").Value = 2132.4
$ws.Range("K122").Value = 6397.200000000001
$ws.Range("M122").Value = -3947.200000000001

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2032
$ws.Range("I132").Value = 1146.0476
$ws.Range("K132").Value = 3438.142800000001
$ws.Range("M132").Value = -908.1428000000005

$ws.Range("H133").Value = 80000
$ws.Range("J133").Value = 80000
$ws.Range("L133").Value = 80000
$ws.Range("N133").Value = -90120
